$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "73.722.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +6.90%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.618.65"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +7.13%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "183.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +12.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "581.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.532"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.196"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +16.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.617.39"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.10%  "
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("E12").Value = "  +7.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "73.625.80"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.94%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000187"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.04%  "
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.076.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +11.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.614.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.94"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +28.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +11.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "371.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +9.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.23"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +14.88%  "
$ws.Range("E23").Value = "  +5.40%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +10.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +11.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.713.39"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0928"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +12.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "515.32"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +19.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +16.72%  "
$ws.Range("E33").Value = "  +5.83%  "
$ws.Range("E34").Value = "  +6.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("E36").Value = "  +12.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "161.10"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.09"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.27"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.37%  "
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.84"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +10.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.324"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.93%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "157.22"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +21.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0873"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +21.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +9.17%  "
$ws.Range("E47").Value = "  +12.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "38.55"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.526"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.14"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +18.79%  "
